$d = $word.ActiveDocument

$replacements = @(
    @("2025-11-01 Saturday", "2025-11-02 Sunday"),
    @("550×9=4950", "811×4=3244"),
    @("261×6=1566", "870×9=7830"),
    @("939×7=6573", "153×4=612"),
    @("679×7=4753", "585×7=4095"),
    @("923×5=4615", "481×9=4329"),
    @("719×3=2157", "711×5=3555"),
    @("987×5=4935", "169×6=1014"),
    @("263×2=526", "330×3=990"),
    @("362×2=724", "305×4=1220"),
    @("553×3=1659", "281×2=562"),
    @("317×8=2536", "204×5=1020"),
    @("176×9=1584", "567×3=1701"),
    @("640×9=5760", "400×6=2400"),
    @("729×7=5103", "666×2=1332"),
    @("966×5=4830", "382×2=764"),
    @("726×3=2178", "605×7=4235"),
    @("154×7=1078", "105×2=210"),
    @("830×2=1660", "598×6=3588"),
    @("851×7=5957", "300×3=900"),
    @("358×3=1074", "834×5=4170"),
    @("849×2=1698", "745×2=1490"),
    @("238×3=714", "611×3=1833"),
    @("935×2=1870", "637×4=2548"),
    @("271×5=1355", "564×4=2256"),
    @("328×7=2296", "965×5=4825")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
